$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C column) date for all existing data rows (2..375)
for ($i = 2; $i -le 375; $i++) {
    $ws.Cells.Item($i, 3).Value = 45184
}

# Row 375 gains an explicit row height (15, custom)
$ws.Rows.Item(375).RowHeight = 15

# Add new row 376 with the new record
$ws.Range("A376").Value = "A 42830-2023"

$ws.Range("B376").Value = 45182
$ws.Range("B376").NumberFormat = "YYYY-MM-DD"

$ws.Range("C376").Value = 45184
$ws.Range("C376").NumberFormat = "YYYY-MM-DD"

$ws.Range("D376").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E376").Value = "VAGGERYD"

$ws.Range("G376").Value = 1.3
$ws.Range("H376").Value = 0
$ws.Range("I376").Value = 0
$ws.Range("J376").Value = 0
$ws.Range("K376").Value = 0
$ws.Range("L376").Value = 0
$ws.Range("M376").Value = 0
$ws.Range("N376").Value = 0
$ws.Range("O376").Value = 0
$ws.Range("P376").Value = 0
$ws.Range("Q376").Value = 0

$ws.Range("R376").Value = ""
$ws.Range("R376").WrapText = $true
